$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above the "004329030" (DANIELA) row, i.e. row 13,
#    shifting all following rows down by one.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row with the account that moved up from
# the bottom of the sheet (now carrying an updated balance). Force the
# account-number column to text so the leading zeros survive, matching
# the "Conta" column elsewhere on the sheet, then drop the temporary
# number-format override so the cell keeps the same (default) styling
# as its neighbours.
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "004381180"
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(13, 2).Value = "HFR"
$ws.Cells.Item(13, 3).Value = 2154.11

# 2. Update BRUNO's (005171652) balance from -19.58 to 0.01.
$ws.Cells.Item(231, 3).Value = 0.01

# 3. Remove the old HFR (004381180 / -34594.91) row further down the
#    sheet - after the insert above it now lives one row lower, at 233.
$ws.Rows.Item(233).Delete()
